$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild the data rows (A2:E40) to reflect the new cleaning step
# ("Drop Rows where prccd is null") inserted into the Stocks pipeline,
# which shifts every subsequent row down by one logical position.
$data = New-Object "object[,]" 39,5
$data[0,0] = 0
$data[0,1] = 'Fundamentals'
$data[0,2] = 'Dropped obvious low value columns -- 1 -- Dates/Codes/Addresses/Etc.'
$data[0,3] = '[''datadate'', ''apdedate'', ''fdate'', ''pdate'', ''fyr'', ''add1'', ''addzip'', ''busdesc'', ''city'', ''conml'', ''ein'', ''fax'', ''fyrc'', ''incorp'', ''loc'', ''phone'', ''state'', ''weburl'', ''exchg'', ''naicsh'', ''sich'', ''naics'', ''sic'', ''stko'']'
$data[0,4] = '(2323, 1745)'
$data[1,0] = 1
$data[1,1] = 'Fundamentals'
$data[1,2] = 'Initial Null Drop -- columns w/ 63% Null or more.'
$data[1,3] = '[''acctchg'', ''acqmeth'', ''adrr'', ''bspr'', ''compst'', ''curuscn'', ''ltcm'', ''ogm'', ''stalt'', ''udpl'', ''acco'', ''acoxar'', ''acqao'', ''acqcshi'', ''acqgdwl'', ''acqic'', ''acqintan'', ''acqinvt'', ''acqlntal'', ''acqniintc'', ''acqppe'', ''acqsc'', ''adpac'', ''aedi'', ''afudcc'', ''afudci'', ''amc'', ''amdc'', ''amgw'', ''ano'', ''aol2'', ''apalch'', ''apb'', ''apc'', ''apofs'', ''aqa'', ''aqd'', ''aqeps'', ''aqi'', ''aqp'', ''aqpl1'', ''aqs'', ''arb'', ''arc'', ''arce'', ''arced'', ''arceeps'', ''artfs'', ''aul3'', ''autxr'', ''balr'', ''banlr'', ''bast'', ''bastr'', ''batr'', ''bcef'', ''bclr'', ''bcltbl'', ''bcnlr'', ''bcrbl'', ''bct'', ''bctbl'', ''bctr'', ''bltbl'', ''ca'', ''capr1'', ''capr2'', ''capr3'', ''cb'', ''cbi'', ''cdpac'', ''cdvc'', ''cfbd'', ''cfere'', ''cfo'', ''cfpdo'', ''cga'', ''cgri'', ''cgti'', ''cgui'', ''chs'', ''clfc'', ''clfx'', ''clg'', ''clis'', ''cll'', ''cllc'', ''clo'', ''clrll'', ''clt'', ''cmp'', ''cnltbl'', ''cpcbl'', ''cpdoi'', ''cpnli'', ''cppbl'', ''cprei'', ''crv'', ''crvnli'', ''cshrc'', ''cshrp'', ''cshrso'', ''cshrt'', ''cshrw'', ''dbi'', ''dcs'', ''depc'', ''dfpac'', ''dfs'', ''dfxa'', ''dlcch'', ''dltsub'', ''dpacb'', ''dpacc'', ''dpacli'', ''dpacls'', ''dpacme'', ''dpacnr'', ''dpaco'', ''dpacre'', ''dpdc'', ''dpltb'', ''dpret'', ''dpsc'', ''dpstb'', ''dptb'', ''dptc'', ''dptic'', ''dpvio'', ''dpvir'', ''drci'', ''dtea'', ''dted'', ''dteeps'', ''dtep'', ''dvdnp'', ''dvpd'', ''dvpdp'', ''dvpibb'', ''dvrpiv'', ''dvrre'', ''dvsco'', ''ea'', ''eiea'', ''emol'', ''esopdlt'', ''excadj'', ''fatd'', ''fatl'', ''fca'', ''fdfr'', ''fea'', ''fel'', ''ffo'', ''ffs'', ''fopt'', ''fsrco'', ''fsrct'', ''fuseo'', ''fuset'', ''gbbl'', ''gdwlam'', ''gdwlia'', ''gdwlid'', ''gdwlieps'', ''gdwlip'', ''geqrv'', ''gla'', ''glcea'', ''glced'', ''glceeps'', ''glcep'', ''gld'', ''gleps'', ''glp'', ''govgr'', ''govtown'', ''gphbl'', ''gplbl'', ''gpobl'', ''gprbl'', ''gptbl'', ''gwo'', ''hedgegl'', ''iaeq'', ''iaeqci'', ''iaeqmi'', ''iafici'', ''iafxi'', ''iafxmi'', ''iali'', ''ialoi'', ''ialti'', ''iamli'', ''iaoi'', ''iapli'', ''iarei'', ''iasci'', ''iasmi'', ''iassi'', ''iasti'', ''iatci'', ''iati'', ''iatmi'', ''iaui'', ''ibbl'', ''ibki'', ''idiis'', ''idilb'', ''idilc'', ''idis'', ''idist'', ''idit'', ''idits'', ''iire'', ''initb'', ''invofs'', ''invreh'', ''invrei'', ''invres'', ''invrm'', ''invwip'', ''iobd'', ''ioi'', ''iore'', ''ip'', ''ipabl'', ''ipc'', ''iphbl'', ''iplbl'', ''ipobl'', ''iptbl'', ''ipti'', ''ipv'', ''irei'', ''irent'', ''irii'', ''irli'', ''irnli'', ''irsi'', ''iseq'', ''iseqc'', ''iseqm'', ''isfi'', ''isfxc'', ''isfxm'', ''isgr'', ''isgt'', ''isgu'', ''islg'', ''islgc'', ''islgm'', ''islt'', ''isng'', ''isngc'', ''isngm'', ''isotc'', ''isoth'', ''isotm'', ''issc'', ''issm'', ''issu'', ''ist'', ''istc'', ''istm'', ''isut'', ''itcc'', ''ivgod'', ''ivi'', ''ivpt'', ''lcabg'', ''lcacl'', ''lcacr'', ''lcag'', ''lcal'', ''lcalt'', ''lcam'', ''lcao'', ''lcast'', ''lcat'', ''lcoxar'', ''lcuacu'', ''li'', ''lif'', ''lloml'', ''lloo'', ''llot'', ''llrci'', ''llrcr'', ''llwoci'', ''llwocr'', ''lno'', ''lol2'', ''lqpl1'', ''lrv'', ''ls'', ''lst'', ''lul3'', ''msvrv'', ''mtl'', ''nat'', ''nco'', ''nfsr'', ''nieci'', ''niint'', ''niintpfc'', ''niintpfp'', ''niit'', ''nim'', ''nio'', ''nipfc'', ''nipfp'', ''nit'', ''nits'', ''npanl'', ''npaore'', ''nparl'', ''npat'', ''nrtxt'', ''nrtxtd'', ''nrtxteps'', ''ob'', ''opili'', ''opincar'', ''opini'', ''opioi'', ''opiri'', ''opiti'', ''optdr'', ''optfvgr'', ''optlife'', ''optprcby'', ''optprcca'', ''optprcex'', ''optprcey'', ''optprcgr'', ''optprcwa'', ''optrfr'', ''optvol'', ''palr'', ''panlr'', ''patr'', ''pcl'', ''pclr'', ''pcnlr'', ''pctr'', ''pdvc'', ''pidom'', ''pifo'', ''pll'', ''pltbl'', ''pncia'', ''pncid'', ''pncieps'', ''pncip'', ''pncwia'', ''pncwid'', ''pncwieps'', ''pncwip'', ''pnlbl'', ''pnli'', ''pobl'', ''ppcbl'', ''ppenb'', ''ppenc'', ''ppenli'', ''ppenls'', ''ppenme'', ''ppennr'', ''ppeno'', ''ppevbb'', ''ppevo'', ''ppevr'', ''pppabl'', ''ppphbl'', ''pppobl'', ''ppptbl'', ''prc'', ''prebl'', ''pri'', ''prodv'', ''prstkcc'', ''prstkpc'', ''prvt'', ''ptbl'', ''ptran'', ''pvcl'', ''pvo'', ''pvon'', ''pvpl'', ''pvt'', ''pwoi'', ''radp'', ''ragr'', ''rari'', ''rati'', ''rca'', ''rcd'', ''rceps'', ''rcl'', ''rcp'', ''rdp'', ''recub'', ''ret'', ''reunr'', ''ris'', ''rll'', ''rlo'', ''rlp'', ''rlri'', ''rlt'', ''rmum'', ''rpag'', ''rra'', ''rrd'', ''rreps'', ''rrp'', ''rvbci'', ''rvbpi'', ''rvbti'', ''rvdo'', ''rvdt'', ''rveqt'', ''rvlrv'', ''rvno'', ''rvnt'', ''rvri'', ''rvsi'', ''rvti'', ''rvtxr'', ''rvupi'', ''rvutx'', ''saa'', ''sal'', ''salepfc'', ''salepfp'', ''sbdc'', ''sc'', ''sco'', ''scstkc'', ''secu'', ''seta'', ''setd'', ''seteps'', ''setp'', ''spid'', ''spieps'', ''spioa'', ''spiop'', ''spstkc'', ''sret'', ''srt'', ''ssnp'', ''stbo'', ''stio'', ''stkcpa'', ''tdc'', ''tdscd'', ''tdsce'', ''tdsg'', ''tdslg'', ''tdsmm'', ''tdsng'', ''tdso'', ''tdss'', ''tdst'', ''tf'', ''tfva'', ''tfvce'', ''tfvl'', ''tie'', ''tii'', ''tlcf'', ''transa'', ''tsa'', ''tsafc'', ''tso'', ''tstkme'', ''txach'', ''txeqa'', ''txeqii'', ''txtubadjust'', ''txtubbegin'', ''txtubend'', ''txtubmax'', ''txtubmin'', ''txtubposdec'', ''txtubposinc'', ''txtubpospdec'', ''txtubpospinc'', ''txtubsettle'', ''txtubsoflimit'', ''txtubtxtr'', ''txtubxintbs'', ''txtubxintis'', ''txva'', ''uaoloch'', ''uaox'', ''uapt'', ''ucaps'', ''uccons'', ''uceq'', ''ucustad'', ''udcopres'', ''udd'', ''udfcc'', ''udmb'', ''udolt'', ''udpco'', ''udpfa'', ''udvp'', ''ufretsd'', ''ugi'', ''ui'', ''uinvt'', ''ulcm'', ''ulco'', ''uniami'', ''unl'', ''unnp'', ''unnpl'', ''unopinc'', ''unwcc'', ''uois'', ''uopi'', ''uopres'', ''updvp'', ''upmcstk'', ''upmpf'', ''upmpfs'', ''upmsubp'', ''upstk'', ''upstkc'', ''upstksf'', ''urect'', ''urectr'', ''urevub'', ''uspi'', ''ustdnc'', ''usubdvp'', ''usubpstk'', ''utfdoc'', ''utfosc'', ''utme'', ''utxfed'', ''uwkcapc'', ''uxinst'', ''uxintd'', ''vpac'', ''vpo'', ''wcapc'', ''wcapch'', ''wda'', ''wdd'', ''wdeps'', ''wdp'', ''xad'', ''xago'', ''xagt'', ''xcom'', ''xcomi'', ''xdepl'', ''xdp'', ''xdvre'', ''xeqo'', ''xindb'', ''xindc'', ''xins'', ''xinst'', ''xintd'', ''xivi'', ''xivre'', ''xlr'', ''xnbi'', ''xnf'', ''xnins'', ''xnitb'', ''xobd'', ''xoi'', ''xoprar'', ''xore'', ''xpp'', ''xpr'', ''xrd'', ''xs'', ''xstf'', ''xstfo'', ''xstfws'', ''xt'', ''xuw'', ''xuwli'', ''xuwnli'', ''xuwoi'', ''xuwrei'', ''xuwti'', ''acchg_fn'', ''afudcc_fn'', ''amc_fn'', ''am_fn'', ''aocidergl_fn'', ''aocipen_fn'', ''aocisecgl_fn'', ''ap_fn'', ''aqa_fn'', ''aqd_fn'', ''aqeps_fn'', ''aqp_fn'', ''aqs_fn'', ''arced_fn'', ''arceeps_fn'', ''arce_fn'', ''artfs_fn'', ''at_fn'', ''capr1_fn'', ''capr3_fn'', ''caps_fn'', ''capx_fn'', ''cdvc_fn'', ''clg_fn'', ''cogs_fn'', ''cshfd_fn'', ''csho_fn'', ''cshpri_fn'', ''cstk_fn'', ''dclo_fn'', ''dcs_fn'', ''dcvsr_fn'', ''dcvsub_fn'', ''dcvt_fn'', ''dd1_fn'', ''dd2_fn'', ''dd3_fn'', ''dd4_fn'', ''dd5_fn'', ''dd_fn'', ''depc_fn'', ''dfs_fn'', ''dlto_fn'', ''dltp_fn'', ''dltt_fn'', ''dn_fn'', ''dpc_fn'', ''dptc_fn'', ''dptic_fn'', ''dp_fn'', ''ds_fn'', ''dtea_fn'', ''dted_fn'', ''dteeps_fn'', ''dtep_fn'', ''dvc_fn'', ''dvp_fn'', ''dvrre_fn'', ''emp_fn'', ''epspi_fn'', ''epspx_fn'', ''esopt_fn'', ''esub_fn'', ''fca_fn'', ''gdwlam_fn'', ''gdwlia_fn'', ''gdwlid_fn'', ''gdwlieps_fn'', ''gdwlip_fn'', ''gla_fn'', ''glcea_fn'', ''glced_fn'', ''glceeps_fn'', ''glcep_fn'', ''gld_fn'', ''gleps_fn'', ''glp_fn'', ''hedgegl_fn'', ''iarei_fn'', ''iati_fn'', ''ibadj_fn'', ''ib_fn'', ''icapt_fn'', ''idiis_fn'', ''idilc_fn'', ''idist_fn'', ''idits_fn'', ''idit_fn'', ''iire_fn'', ''invt_fn'', ''ip_fn'', ''ist_fn'', ''itcc_fn'', ''itci_fn'', ''ivi_fn'', ''lcao_fn'', ''lcuacu_fn'', ''lse_fn'', ''mib_fn'', ''mrc1_fn'', ''mrc2_fn'', ''mrc3_fn'', ''mrc4_fn'', ''mrc5_fn'', ''msa_fn'', ''ni_fn'', ''nrtxtd_fn'', ''nrtxteps_fn'', ''nrtxt_fn'', ''oancf_fn'', ''opeps_fn'', ''oprepsx_fn'', ''optdr_fn'', ''optfvgr_fn'', ''optlife_fn'', ''optprcby_fn'', ''optprcca_fn'', ''optprcex_fn'', ''optprcey_fn'', ''optprcgr_fn'', ''optprcwa_fn'', ''optrfr_fn'', ''optvol_fn'', ''pcl_fn'', ''pdvc_fn'', ''pidom_fn'', ''pifo_fn'', ''pncad_fn'', ''pncaeps_fn'', ''pnca_fn'', ''pncia_fn'', ''pncid_fn'', ''pncieps_fn'', ''pncwia_fn'', ''pncwid_fn'', ''pncwieps_fn'', ''ppegt_fn'', ''ppent_fn'', ''ppeveb_fn'', ''prcad_fn'', ''prcaeps_fn'', ''prca_fn'', ''pstkl_fn'', ''pstkrv_fn'', ''pvcl_fn'', ''ragr_fn'', ''rca_fn'', ''rcd_fn'', ''rceps_fn'', ''rcl_fn'', ''rcp_fn'', ''rdip_fn'', ''recta_fn'', ''reuna_fn'', ''revt_fn'', ''re_fn'', ''rra_fn'', ''rrd_fn'', ''rreps_fn'', ''rrp_fn'', ''sale_fn'', ''secu_fn'', ''seta_fn'', ''setd_fn'', ''seteps_fn'', ''setp_fn'', ''spced_fn'', ''spceeps_fn'', ''spce_fn'', ''spid_fn'', ''spieps_fn'', ''spioa_fn'', ''spiop_fn'', ''stkcpa_fn'', ''tdc_fn'', ''tdst_fn'', ''tie_fn'', ''transa_fn'', ''tstkc_fn'', ''tstk_fn'', ''txeqii_fn'', ''txfed_fn'', ''txfo_fn'', ''txo_fn'', ''txp_fn'', ''txs_fn'', ''txt_fn'', ''uapt_fn'', ''ucaps_fn'', ''ui_fn'', ''unwcc_fn'', ''vpac_fn'', ''vpo_fn'', ''wda_fn'', ''wdd_fn'', ''wdeps_fn'', ''wdp_fn'', ''xintopt_fn'', ''xint_fn'', ''xlr_fn'', ''xoptd_fn'', ''xopteps_fn'', ''xpr_fn'', ''xrd_fn'', ''xrent_fn'', ''xsga_fn'', ''xt_fn'', ''acchg_dc'', ''acco_dc'', ''acdo_dc'', ''aco_dc'', ''acominc_dc'', ''acox_dc'', ''acqao_dc'', ''acqgdwl_dc'', ''acqic_dc'', ''acqintan_dc'', ''acqinvt_dc'', ''acqlntal_dc'', ''acqniintc_dc'', ''acqppe_dc'', ''acqsc_dc'', ''adpac_dc'', ''aedi_dc'', ''aldo_dc'', ''am_dc'', ''amc_dc'', ''amgw_dc'', ''ano_dc'', ''ao_dc'', ''aocidergl_dc'', ''aocipen_dc'', ''aocisecgl_dc'', ''aol2_dc'', ''aoloch_dc'', ''aox_dc'', ''ap_dc'', ''apalch_dc'', ''apofs_dc'', ''aqa_dc'', ''aqc_dc'', ''aqd_dc'', ''aqeps_dc'', ''aqi_dc'', ''aqp_dc'', ''aqpl1_dc'', ''aqs_dc'', ''arc_dc'', ''artfs_dc'', ''aul3_dc'', ''balr_dc'', ''banlr_dc'', ''bast_dc'', ''bastr_dc'', ''batr_dc'', ''bcef_dc'', ''bclr_dc'', ''bcnlr_dc'', ''bcrbl_dc'', ''bct_dc'', ''bctbl_dc'', ''bctr_dc'', ''bkvlps_dc'', ''bltbl_dc'', ''ca_dc'', ''caps_dc'', ''capx_dc'', ''capxv_dc'', ''cb_dc'', ''cbi_dc'', ''cdvc_dc'', ''ceq_dc'', ''ceql_dc'', ''ceqt_dc'', ''cfbd_dc'', ''cfere_dc'', ''cfo_dc'', ''cfpdo_dc'', ''cgti_dc'', ''ch_dc'', ''che_dc'', ''chech_dc'', ''chs_dc'', ''cicurr_dc'', ''cidergl_dc'', ''cipen_dc'', ''cisecgl_dc'', ''cld2_dc'', ''cld3_dc'', ''cld4_dc'', ''cld5_dc'', ''clfc_dc'', ''clfx_dc'', ''clg_dc'', ''clis_dc'', ''cll_dc'', ''cllc_dc'', ''cmp_dc'', ''cnltbl_dc'', ''cpnli_dc'', ''cprei_dc'', ''crvnli_dc'', ''cshfd_dc'', ''cshi_dc'', ''csho_dc'', ''cshpri_dc'', ''cshr_dc'', ''cshrc_dc'', ''cshrt_dc'', ''cshrw_dc'', ''cstk_dc'', ''cstkcv_dc'', ''dbi_dc'', ''dc_dc'', ''dclo_dc'', ''dcom_dc'', ''dcpstk_dc'', ''dcs_dc'', ''dcvsr_dc'', ''dcvsub_dc'', ''dcvt_dc'', ''dd_dc'', ''dd1_dc'', ''dd2_dc'', ''dd3_dc'', ''dd4_dc'', ''dd5_dc'', ''depc_dc'', ''dfpac_dc'', ''dfs_dc'', ''dfxa_dc'', ''dlc_dc'', ''dlcch_dc'', ''dltis_dc'', ''dlto_dc'', ''dltp_dc'', ''dltr_dc'', ''dltsub_dc'', ''dltt_dc'', ''dm_dc'', ''dn_dc'', ''do_dc'', ''donr_dc'', ''dp_dc'', ''dpacb_dc'', ''dpacc_dc'', ''dpacls_dc'', ''dpacme_dc'', ''dpacnr_dc'', ''dpact_dc'', ''dpc_dc'', ''dpdc_dc'', ''dpsc_dc'', ''dptb_dc'', ''dptc_dc'', ''dptic_dc'', ''dpvieb_dc'', ''dpvio_dc'', ''dpvir_dc'', ''drc_dc'', ''drlt_dc'', ''ds_dc'', ''dtea_dc'', ''dted_dc'', ''dteeps_dc'', ''dtep_dc'', ''dudd_dc'', ''dv_dc'', ''dvc_dc'', ''dvp_dc'', ''dvpa_dc'', ''dvpdp_dc'', ''dvpibb_dc'', ''dvrre_dc'', ''dxd2_dc'', ''dxd3_dc'', ''dxd4_dc'', ''dxd5_dc'', ''ea_dc'', ''ebitda_dc'', ''eiea_dc'', ''emp_dc'', ''epsfi_dc'', ''epsfx_dc'', ''epspi_dc'', ''epspx_dc'', ''esopct_dc'', ''esub_dc'', ''esubc_dc'', ''excadj_dc'', ''exre_dc'', ''fatb_dc'', ''fatc_dc'', ''fate_dc'', ''fatl_dc'', ''fatn_dc'', ''fato_dc'', ''fatp_dc'', ''fca_dc'', ''fdfr_dc'', ''ffs_dc'', ''fopo_dc'', ''fsrco_dc'', ''fuseo_dc'', ''gdwl_dc'', ''gdwlia_dc'', ''gdwlid_dc'', ''gdwlieps_dc'', ''gdwlip_dc'', ''gla_dc'', ''glcea_dc'', ''glcep_dc'', ''gld_dc'', ''gleps_dc'', ''glp_dc'', ''govgr_dc'', ''gphbl_dc'', ''gplbl_dc'', ''gprbl_dc'', ''hedgegl_dc'', ''iaeq_dc'', ''iaeqci_dc'', ''iaeqmi_dc'', ''ialti_dc'', ''iamli_dc'', ''iapli_dc'', ''iarei_dc'', ''iati_dc'', ''ib_dc'', ''ibadj_dc'', ''ibcom_dc'', ''ibki_dc'', ''icapt_dc'', ''idiis_dc'', ''idilb_dc'', ''idilc_dc'', ''idis_dc'', ''idist_dc'', ''idit_dc'', ''idits_dc'', ''iire_dc'', ''initb_dc'', ''intan_dc'', ''intano_dc'', ''intc_dc'', ''intpn_dc'', ''invch_dc'', ''invfg_dc'', ''invreh_dc'', ''invrei_dc'', ''invres_dc'', ''invrm_dc'', ''invt_dc'', ''invwip_dc'', ''iobd_dc'', ''ioi_dc'', ''iore_dc'', ''ip_dc'', ''ipabl_dc'', ''iphbl_dc'', ''iplbl_dc'', ''ipti_dc'', ''ipv_dc'', ''irei_dc'', ''irent_dc'', ''irli_dc'', ''irnli_dc'', ''irsi_dc'', ''iseq_dc'', ''iseqc_dc'', ''iseqm_dc'', ''isfi_dc'', ''isfxc_dc'', ''isfxm_dc'', ''isgt_dc'', ''islg_dc'', ''islgc_dc'', ''islgm_dc'', ''isng_dc'', ''isngc_dc'', ''isngm_dc'', ''issc_dc'', ''issm_dc'', ''issu_dc'', ''ist_dc'', ''itcb_dc'', ''itcc_dc'', ''itci_dc'', ''ivaeq_dc'', ''ivao_dc'', ''ivch_dc'', ''ivi_dc'', ''ivpt_dc'', ''ivst_dc'', ''ivstch_dc'', ''lcabg_dc'', ''lcacl_dc'', ''lcacr_dc'', ''lcag_dc'', ''lcal_dc'', ''lcam_dc'', ''lcat_dc'', ''lco_dc'', ''lcox_dc'', ''lct_dc'', ''lcuacu_dc'', ''lifr_dc'', ''llrcr_dc'', ''llwocr_dc'', ''lno_dc'', ''lo_dc'', ''lol2_dc'', ''lqpl1_dc'', ''lt_dc'', ''lul3_dc'', ''mib_dc'', ''mibn_dc'', ''mibt_dc'', ''mii_dc'', ''mrc1_dc'', ''mrc2_dc'', ''mrc3_dc'', ''mrc4_dc'', ''mrc5_dc'', ''mrct_dc'', ''mrcta_dc'', ''msa_dc'', ''nco_dc'', ''ni_dc'', ''nieci_dc'', ''niintpfc_dc'', ''niintpfp_dc'', ''nio_dc'', ''nipfc_dc'', ''nipfp_dc'', ''nit_dc'', ''nopi_dc'', ''nopio_dc'', ''np_dc'', ''npaore_dc'', ''nparl_dc'', ''nrtxt_dc'', ''nrtxtd_dc'', ''nrtxteps_dc'', ''ob_dc'', ''oiadp_dc'', ''oibdp_dc'', ''opili_dc'', ''opini_dc'', ''opioi_dc'', ''opiri_dc'', ''optca_dc'', ''optexd_dc'', ''optgr_dc'', ''optosey_dc'', ''optprcca_dc'', ''optprcex_dc'', ''optprcgr_dc'', ''palr_dc'', ''panlr_dc'', ''patr_dc'', ''pcl_dc'', ''pclr_dc'', ''pcnlr_dc'', ''pctr_dc'', ''pdvc_dc'', ''pi_dc'', ''pidom_dc'', ''pifo_dc'', ''pll_dc'', ''pltbl_dc'', ''pnlbl_dc'', ''pnrsho_dc'', ''pobl_dc'', ''ppegt_dc'', ''ppenb_dc'', ''ppenc_dc'', ''ppenli_dc'', ''ppenls_dc'', ''ppenme_dc'', ''ppennr_dc'', ''ppeno_dc'', ''ppent_dc'', ''ppevbb_dc'', ''ppeveb_dc'', ''ppevo_dc'', ''ppevr_dc'', ''ppphbl_dc'', ''prebl_dc'', ''pri_dc'', ''prodv_dc'', ''prsho_dc'', ''prstkc_dc'', ''prstkcc_dc'', ''prstkpc_dc'', ''pstk_dc'', ''pstkc_dc'', ''pstkl_dc'', ''pstkn_dc'', ''pstkr_dc'', ''pstkrv_dc'', ''ptran_dc'', ''pvon_dc'', ''pvpl_dc'', ''pvt_dc'', ''pwoi_dc'', ''ragr_dc'', ''rati_dc'', ''rca_dc'', ''rcd_dc'', ''rceps_dc'', ''rcl_dc'', ''rcp_dc'', ''rdip_dc'', ''rdipa_dc'', ''rdipd_dc'', ''re_dc'', ''rea_dc'', ''reajo_dc'', ''recch_dc'', ''recco_dc'', ''recd_dc'', ''rect_dc'', ''recta_dc'', ''rectr_dc'', ''reunr_dc'', ''ris_dc'', ''rll_dc'', ''rlri_dc'', ''rlt_dc'', ''rpag_dc'', ''rvbci_dc'', ''rvbti_dc'', ''rvdt_dc'', ''rveqt_dc'', ''rvri_dc'', ''rvsi_dc'', ''rvti_dc'', ''rvupi_dc'', ''saa_dc'', ''sal_dc'', ''salepfc_dc'', ''salepfp_dc'', ''sbdc_dc'', ''sc_dc'', ''scstkc_dc'', ''secu_dc'', ''seq_dc'', ''seta_dc'', ''setd_dc'', ''seteps_dc'', ''setp_dc'', ''siv_dc'', ''spi_dc'', ''spid_dc'', ''spieps_dc'', ''spioa_dc'', ''spiop_dc'', ''sppe_dc'', ''sppiv_dc'', ''spstkc_dc'', ''ssnp_dc'', ''sstk_dc'', ''stbo_dc'', ''stio_dc'', ''stkco_dc'', ''stkcpa_dc'', ''tdc_dc'', ''tdscd_dc'', ''tdsce_dc'', ''tdsg_dc'', ''tdslg_dc'', ''tdsmm_dc'', ''tdsng_dc'', ''tdst_dc'', ''tf_dc'', ''tfva_dc'', ''tfvce_dc'', ''tfvl_dc'', ''tlcf_dc'', ''transa_dc'', ''tstk_dc'', ''tstkc_dc'', ''tstkme_dc'', ''tstkn_dc'', ''tstkp_dc'', ''txach_dc'', ''txbco_dc'', ''txbcof_dc'', ''txc_dc'', ''txdb_dc'', ''txdba_dc'', ''txdbca_dc'', ''txdbcl_dc'', ''txdc_dc'', ''txdfed_dc'', ''txdfo_dc'', ''txdi_dc'', ''txditc_dc'', ''txds_dc'', ''txeqa_dc'', ''txfed_dc'', ''txfo_dc'', ''txndb_dc'', ''txndba_dc'', ''txndbl_dc'', ''txndbr_dc'', ''txp_dc'', ''txpd_dc'', ''txr_dc'', ''txs_dc'', ''txt_dc'', ''txtubadjust_dc'', ''txtubbegin_dc'', ''txtubend_dc'', ''txtubmax_dc'', ''txtubmin_dc'', ''txtubposdec_dc'', ''txtubposinc_dc'', ''txtubpospdec_dc'', ''txtubpospinc_dc'', ''txtubsettle_dc'', ''txtubsoflimit_dc'', ''txtubtxtr_dc'', ''txtubxintbs_dc'', ''txtubxintis_dc'', ''txw_dc'', ''uaoloch_dc'', ''uaox_dc'', ''uapt_dc'', ''ucaps_dc'', ''uccons_dc'', ''ucustad_dc'', ''udd_dc'', ''udfcc_dc'', ''udmb_dc'', ''udpco_dc'', ''udvp_dc'', ''ui_dc'', ''uinvt_dc'', ''unopinc_dc'', ''unwcc_dc'', ''uopi_dc'', ''uopres_dc'', ''updvp_dc'', ''upmcstk_dc'', ''upmpf_dc'', ''upmpfs_dc'', ''upmsubp_dc'', ''upstk_dc'', ''upstkc_dc'', ''upstksf_dc'', ''urect_dc'', ''urectr_dc'', ''urevub_dc'', ''uspi_dc'', ''ustdnc_dc'', ''usubdvp_dc'', ''usubpstk_dc'', ''utme_dc'', ''utxfed_dc'', ''uwkcapc_dc'', ''uxinst_dc'', ''uxintd_dc'', ''vpac_dc'', ''vpo_dc'', ''wcapc_dc'', ''wda_dc'', ''wdd_dc'', ''wdeps_dc'', ''wdp_dc'', ''xacc_dc'', ''xad_dc'', ''xago_dc'', ''xagt_dc'', ''xcom_dc'', ''xcomi_dc'', ''xdepl_dc'', ''xdp_dc'', ''xdvre_dc'', ''xeqo_dc'', ''xi_dc'', ''xido_dc'', ''xidoc_dc'', ''xindb_dc'', ''xindc_dc'', ''xins_dc'', ''xinst_dc'', ''xint_dc'', ''xintd_dc'', ''xintopt_dc'', ''xivi_dc'', ''xivre_dc'', ''xlr_dc'', ''xnins_dc'', ''xnitb_dc'', ''xobd_dc'', ''xoi_dc'', ''xoptd_dc'', ''xopteps_dc'', ''xore_dc'', ''xpp_dc'', ''xpr_dc'', ''xrd_dc'', ''xrent_dc'', ''xs_dc'', ''xsga_dc'', ''xstf_dc'', ''xstfo_dc'', ''xuw_dc'', ''xuwli_dc'', ''xuwnli_dc'', ''xuwoi_dc'', ''xuwrei_dc'', ''add2'', ''add3'', ''add4'', ''county'', ''dlrsn'', ''prican'', ''prirow'', ''dldte'', ''ipodate'']'
$data[1,4] = '(2323, 327)'
$data[2,0] = 2
$data[2,1] = 'Fundamentals'
$data[2,2] = 'Dropped columns with only 1 unique value '
$data[2,3] = '[''indfmt'', ''consol'', ''popsrc'', ''curcd'', ''final'', ''scf'', ''acchg'', ''aocisecgl'', ''esopr'', ''itcb'', ''xoptd'', ''xopteps'', ''rank'', ''rdipa_fn'', ''stkco_fn'', ''gsector'']'
$data[2,4] = '(2323, 311)'
$data[3,0] = 3
$data[3,1] = 'Fundamentals'
$data[3,2] = 'Dropping all that are not standard statements -- value of others captured in restatement summary'
$data[3,3] = 'applies to all cols, roughly 50% of rows dropped'
$data[3,4] = '(1243, 311)'
$data[4,0] = 4
$data[4,1] = 'Fundamentals'
$data[4,2] = 'Second Null Drop -- columns w/ 20% Null or more.'
$data[4,3] = '[''cld2'', ''cld3'', ''cld4'', ''cld5'', ''cshr'', ''dd3'', ''dd4'', ''dd5'', ''dltp'', ''dxd2'', ''dxd3'', ''dxd4'', ''dxd5'', ''fatb'', ''fate'', ''fato'', ''fatp'', ''invfg'', ''invo'', ''itci'', ''mrc2'', ''mrc3'', ''mrc4'', ''mrc5'', ''mrct'', ''optca'', ''optex'', ''optexd'', ''optgr'', ''optosby'', ''optosey'', ''recd'', ''sppe'', ''stkco'', ''txc'', ''txdfed'', ''txdfo'', ''txds'', ''txfed'', ''txs'', ''xrent'', ''mkvalt'', ''spcindcd'', ''spcseccd'', ''spcsrc'']'
$data[4,4] = '(1243, 266)'
$data[5,0] = 5
$data[5,1] = 'Fundamentals'
$data[5,2] = 'Dropping Columns w/ Only 2 values'
$data[5,3] = '[''datafmt'', ''ismod'', ''upd'', ''rdip'', ''rdipa'', ''rdipd'', ''rdipeps'', ''tstkp'', ''txndbr'', ''xi'', ''xintopt'', ''costat'', ''rdipd_fn'', ''rdipeps_fn'', ''idbflag'']'
$data[5,4] = '(1243, 251)'
$data[6,0] = 6
$data[6,1] = 'Fundamentals'
$data[6,2] = 'Second Manual Column Reduce - Dropping Identified Low Value Columns'
$data[6,3] = '[''Unnamed: 0'', ''conm'', ''au'', ''auopic'', ''pddur'', ''dpact_fn'', ''ggroup'', ''gind'', ''fic'', ''priusa'']'
$data[6,4] = '(1243, 241)'
$data[7,0] = 7
$data[7,1] = 'Fundamentals'
$data[7,2] = 'Third Manual Column Reduce  - Targetted Reduce of Granular Features w/ Total'
$data[7,3] = '[''lco'', ''lcox'', ''lcoxdr'', ''lo'', ''loxdr'', ''mib'', ''mibn'', ''ppeveb'', ''pstkc'', ''pstkl'', ''pstkn'', ''pstkr'', ''tstkc'', ''tstkn'', ''che'', ''cicurr'', ''cidergl'', ''cimii'', ''ciother'', ''cipen'', ''cisecgl'', ''citotal'', ''dpc'', ''dpvieb'', ''dv'', ''dvp'', ''dvpa'', ''ib'', ''epspx'', ''esopnr'', ''esopt'', ''ibadj'', ''ibc'', ''ibcom'', ''ibmii'', ''recch'', ''recco'', ''rectr'', ''reuna'', ''sale'', ''spced'', ''spceeps'', ''cshtr_c'', ''dvpsp_c'', ''dvpsx_c'', ''prcc_c'', ''prch_c'', ''prcl_c'', ''adjex_c'', ''acdo'', ''aco'', ''acodo'', ''acox'', ''aldo'', ''aocidergl'', ''aociother'', ''aocipen'', ''aodo'', ''aox'', ''ap'', ''ceql'', ''intc'', ''ivaco'', ''ivaeq'', ''ivao'', ''ivch'', ''ivst'', ''ivstch'', ''pncad'', ''pncaeps'', ''prcad'', ''prcaeps'', ''xido'', ''xidoc'', ''ajex'', ''ajp'', ''cshfd'', ''cshi'', ''csho'', ''cstk'', ''cstkcv'', ''cstke'', ''dclo'', ''dcom'', ''dcvsr'', ''dcvsub'', ''dcvt'', ''dd'', ''dd1'', ''dd2'', ''dltis'', ''dlto'', ''dm'', ''dn'', ''ds'', ''dudd'', ''fatc'', ''fatc'', ''fatn'', ''fiao'', ''fopox'', ''fyear'', ''intano'', ''mrc1'', ''mrcta'', ''niadj'', ''nopio'', ''oiadp'', ''oibdp'', ''oprepsx'', ''pnrsho'', ''prsho'', ''ppent'', ''pstkrv'', ''txbco'', ''txbcof'', ''txdba'', ''txdbca'', ''txdbcl'', ''txdc'', ''txdi'', ''txditc'', ''txndba'', ''txndb'', ''txndbl'', ''txo'', ''txp'', ''txpd'', ''txr'', ''acctstd'', ''am'', ''capxv'', ''dc'', ''diladj'', ''do'', ''donr'', ''emp'', ''esub'', ''exre'', ''lifr'', ''mibt'', ''prstkc'', ''seqo'', ''spi'', ''cshtr_f'', ''dvpsp_f'', ''dvpsx_f'', ''prcc_f'', ''prch_f'', ''prcl_f'', ''adjex_f'', ''ebitda'', ''acominc'', ''act'', ''ao'', ''aqc'', ''cfoso'']'
$data[7,4] = '(1243, 85)'
$data[8,0] = 8
$data[8,1] = 'Fundamentals'
$data[8,2] = 'Summary Features  - Used to create features, then dropped'
$data[8,3] = '[''drc'', ''drlt'', ''rea'', ''reajo'', ''recta'', ''txdb'', ''txfo'', ''txt'', ''txw'']'
$data[8,4] = '(1243, 78)'
$data[9,0] = 9
$data[9,1] = 'Fundamentals'
$data[9,2] = 'Group by GVKey'
$data[9,3] = 'All'
$data[9,4] = '(348, 147)'
$data[10,0] = 10
$data[10,1] = 'Fundamentals'
$data[10,2] = 'Join Restatement Data'
$data[10,3] = 'Restatement Cols'
$data[10,4] = '(348, 152)'
$data[11,0] = 11
$data[11,1] = 'Ratings'
$data[11,2] = 'Capturing Original Shape'
$data[11,3] = 'N/A'
$data[11,4] = '(123678, 21)'
$data[12,0] = 12
$data[12,1] = 'Ratings'
$data[12,2] = 'Trim Rows by GSector == 30'
$data[12,3] = 'All'
$data[12,4] = '(17902, 21)'
$data[13,0] = 13
$data[13,1] = 'Ratings'
$data[13,2] = 'Dropped obvious low value columns -- 1 -- Dates/Codes/Addresses/Etc.'
$data[13,3] = '[''city'', ''ggroup'', ''gind'', ''gsector'', ''gsubind'', ''loc'', ''naics'', ''state'', ''conm'', ''tic'', ''idbflag'', ''sic'', ''spcindcd'', ''spcseccd'', ''conml'']'
$data[13,4] = '(17902, 6)'
$data[14,0] = 14
$data[14,1] = 'Ratings'
$data[14,2] = 'Dropped Nulls -- 1 -- More than 50% Null -- Drops more granular Ratings'
$data[14,3] = '[''splticrm'', ''spsdrm'', ''spsticrm'']'
$data[14,4] = '(17902, 3)'
$data[15,0] = 15
$data[15,1] = 'Ratings'
$data[15,2] = 'Group by GVKey'
$data[15,3] = '[''city'', ''ggroup'', ''gind'', ''gsector'', ''gsubind'', ''loc'', ''naics'', ''state'', ''conm'', ''tic'', ''idbflag'', ''sic'', ''spcindcd'', ''spcseccd'', ''conml'']'
$data[15,4] = '(452, 2)'
$data[16,0] = 16
$data[16,1] = 'Securities'
$data[16,2] = 'Capturing Original Shape'
$data[16,3] = 'N/A'
$data[16,4] = '(32382, 56)'
$data[17,0] = 17
$data[17,1] = 'Securities'
$data[17,2] = 'Dropped obvious low value columns -- 1 -- Dates/Codes/Addresses/Etc.'
$data[17,3] = '[''Unnamed: 0'', ''conm'', ''city'', ''conml'', ''cyear'', ''ggroup'', ''gind'', ''gsector'', ''gsubind'', ''loc'', ''naics'', ''sic'', ''state'']'
$data[17,4] = '(32382, 43)'
$data[18,0] = 18
$data[18,1] = 'Securities'
$data[18,2] = 'Dropped Nulls -- 1 -- More than 20% Null'
$data[18,3] = '[''isalrt'', ''spgim'', ''spiim'', ''spmim'', ''cheqvm'', ''curcddvm'', ''dvpspm'', ''dvpsxm'', ''dvrate'', ''csfsm'', ''navm'', ''rawpm'', ''rawxm'', ''sph100'', ''sphcusip'', ''sphiid'', ''sphmid'', ''sphname'', ''sphsec'', ''sphtic'', ''sphvg'', ''cshoq'', ''adrrm'']'
$data[18,4] = '(32382, 20)'
$data[19,0] = 19
$data[19,1] = 'Securities'
$data[19,2] = 'Dropping Columns w/ Only 2 values'
$data[19,3] = '[''primiss'', ''curcdm'', ''mkvalincl'', ''costat'']'
$data[19,4] = '(32382, 16)'
$data[20,0] = 20
$data[20,1] = 'Securities'
$data[20,2] = 'Manual Drop -- 2 -- Low value or duplicated by other datasets'
$data[20,3] = '[''iid'', ''tic'', ''cmth'', ''exchg'', ''tpci'', ''cshtrm'', ''prccm'', ''prchm'', ''prclm'', ''cshom'']'
$data[20,4] = '(32382, 6)'
$data[21,0] = 21
$data[21,1] = 'Securities'
$data[21,2] = 'Group by GVKey'
$data[21,3] = '[''iid'', ''tic'', ''cmth'', ''exchg'', ''tpci'', ''cshtrm'', ''prccm'', ''prchm'', ''prclm'', ''cshom'']'
$data[21,4] = '(690, 7)'
$data[22,0] = 22
$data[22,1] = 'Stocks'
$data[22,2] = 'Capturing Original Shape'
$data[22,3] = 'N/A'
$data[22,4] = '(673564, 77)'
$data[23,0] = 23
$data[23,1] = 'Stocks'
$data[23,2] = 'Dropped Nulls -- 1 -- More than 20% Null'
$data[23,3] = '[''curcddv'', ''capgn'', ''cheqv'', ''div'', ''divd'', ''divdpaydateind'', ''divsp'', ''dvrated'', ''paydateind'', ''anncdate'', ''capgnpaydate'', ''cheqvpaydate'', ''divdpaydate'', ''divsppaydate'', ''paydate'', ''recorddate'', ''adrrc'', ''dvi'', ''eps'', ''epsmo'', ''prcod'', ''add2'', ''add3'', ''add4'', ''county'', ''dlrsn'', ''ein'', ''fax'', ''incorp'', ''prican'', ''prirow'', ''spcindcd'', ''spcseccd'', ''spcsrc'', ''state'', ''dldte'', ''ipodate'']'
$data[23,4] = '(673564, 40)'
$data[24,0] = 24
$data[24,1] = 'Stocks'
$data[24,2] = 'Dropped obvious low value columns -- 1 -- Dates/Codes/Addresses/Etc.'
$data[24,3] = '[''Unnamed: 0'', ''tic'', ''conm'', ''curcdd'', ''fic'', ''add1'', ''addzip'', ''busdesc'', ''city'', ''conml'', ''fyrc'', ''ggroup'', ''gind'', ''gsector'', ''gsubind'', ''idbflag'', ''loc'', ''naics'', ''phone'', ''sic'', ''weburl'', ''prcstd'', ''cusip'', ''priusa'']'
$data[24,4] = '(673564, 16)'
$data[25,0] = 25
$data[25,1] = 'Stocks'
$data[25,2] = 'Dropping Columns w/ Only 2 values'
$data[25,3] = '[''secstat'', ''costat'', ''stko'']'
$data[25,4] = '(673564, 13)'
$data[26,0] = 26
$data[26,1] = 'Stocks'
$data[26,2] = 'Manual Drop -- 2 -- Low value or duplicated by other datasets'
$data[26,3] = '[''iid'', ''exchg'', ''tpci'', ''cik'', ''cshtrd'', ''ajexdi'', ''trfd'', ''cshoc'']'
$data[26,4] = '(673564, 5)'
$data[27,0] = 27
$data[27,1] = 'Stocks'
$data[27,2] = 'Drop Rows where prccd is null -- usually an additional row for a dividend. See comments.'
$data[27,3] = 'prccd'
$data[27,4] = '(673301, 5)'
$data[28,0] = 28
$data[28,1] = 'Stocks'
$data[28,2] = 'Group by GVKey'
$data[28,3] = '[''iid'', ''exchg'', ''tpci'', ''cik'', ''cshtrd'', ''ajexdi'', ''trfd'', ''cshoc'']'
$data[28,4] = '(689, 7)'
$data[29,0] = 29
$data[29,1] = 'Stocks'
$data[29,2] = 'Add Calculated Fields -- replace variables no longer needed'
$data[29,3] = '[''iid'', ''exchg'', ''tpci'', ''cik'', ''cshtrd'', ''ajexdi'', ''trfd'', ''cshoc'']'
$data[29,4] = '(689, 6)'
$data[30,0] = 30
$data[30,1] = 'SCA_Filings'
$data[30,2] = 'Capturing Original Shape'
$data[30,3] = 'N/A'
$data[30,4] = '(1892, 6)'
$data[31,0] = 31
$data[31,1] = 'SCA_Filings'
$data[31,2] = 'Manual Drop -- 1 -- Low value or duplicated by other datasets'
$data[31,3] = '[''FilingName'', ''Exchange'', ''FilingYear'', ''Dismissed'']'
$data[31,4] = '(1892, 2)'
$data[32,0] = 32
$data[32,1] = 'SCA_Filings'
$data[32,2] = 'Group by Ticker'
$data[32,3] = '[''FilingName'', ''Exchange'', ''FilingYear'', ''Dismissed'']'
$data[32,4] = '(1506, 2)'
$data[33,0] = 33
$data[33,1] = 'SCA_Filings'
$data[33,2] = 'Create Boolean Response Var'
$data[33,3] = '[''FilingName'', ''Exchange'', ''FilingYear'', ''Dismissed'']'
$data[33,4] = '(1506, 3)'
$data[34,0] = 34
$data[34,1] = 'completeDF'
$data[34,2] = 'Make CompleteDF from a copy of fundaments'
$data[34,3] = 'All fundamentals cols'
$data[34,4] = '(348, 152)'
$data[35,0] = 35
$data[35,1] = 'completeDF'
$data[35,2] = 'Add Stocks to CompleteDF'
$data[35,3] = 'fundamentals + stocks'
$data[35,4] = '(348, 157)'
$data[36,0] = 36
$data[36,1] = 'completeDF'
$data[36,2] = 'Add Securities to CompleteDF'
$data[36,3] = 'fundamentals + stocks + securites'
$data[36,4] = '(348, 163)'
$data[37,0] = 37
$data[37,1] = 'completeDF'
$data[37,2] = 'Add Ratings to CompleteDF'
$data[37,3] = 'fundamentals + stocks + securites + ratings'
$data[37,4] = '(348, 164)'
$data[38,0] = 38
$data[38,1] = 'completeDF'
$data[38,2] = 'Add SCA Filings to CompleteDF'
$data[38,3] = 'fundamentals + stocks + securites + ratings + SCA'
$data[38,4] = '(348, 167)'

$ws.Range("A2:E40").Value = $data

Write-Output "done"